$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The codeforiati:group-code (D) and codeforiati:group-name (E) columns were
# swapped for the header row and every data row (1-94): what used to be in D
# is now in E and vice versa.
for ($r = 1; $r -le 94; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value = $eVal
    $eCell.Value = $dVal
}
